$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.114.79"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "3.787.59"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "3.787.10"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.52"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000281"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "4.420.56"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "3.790.14"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.73"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "67.953.05"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.47"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "3.934.52"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.60"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.29"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "3.743.75"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.76"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.19%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.96"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "402.06"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.19%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000280"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.00%  "
